# Refresh the crypto price/volume snapshot (GitHub Actions scrape update).
# For each changed cell: plain numeric-looking strings (e.g. "1.002") are
# written with a temporary Text number format so Excel keeps them as literal
# strings (matching the source inlineStr cells) instead of silently parsing
# them into doubles / losing trailing zeros; the style is then restored so no
# stray formatting is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '24.580.02'
$ws.Range('E2').Value = '  +3.10%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '1.695.42'
$ws.Range('E3').Value = '  +1.90%  '

# Row 4: TetherUSD
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.12%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.44'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.23%  '

# Row 6: USDC
$ws.Range('E6').Value = '  +0.11%  '

# Row 7: XRP
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3950'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.86%  '

# Row 8: Cardano
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4009'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.18%  '

# Row 9: Polygon
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.520'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.51%  '

# Row 10: BinanceUSD
$ws.Range('E10').Value = '  +0.06%  '

# Row 11: OKB
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.21'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.53%  '

# Row 12: Dogecoin
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08758'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.02%  '

# Row 13: Polkadot
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.225'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.72%  '

# Row 14: Solana
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.23'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.50%  '

# Row 15: Chainlink
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.167'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +11.64%  '

# Row 16: ShibaInu
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001311'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.40%  '

# Row 17: WrappedEther
$ws.Range('D17').Value = '1.694.15'

# Row 18: Litecoin
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '99.83'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.68%  '

# Row 19: TRON
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07068'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.93%  '

# Row 20: Avalanche
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.63'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.93%  '

# Row 21: Uniswap
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.024'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.02%  '

# Row 22: Dai
$ws.Range('E22').Value = '  -0.22%  '

# Row 23: Cosmos
$ws.Range('E23').Value = '  +3.00%  '

# Row 24: WrappedBTC
$ws.Range('D24').Value = '24.584.38'
$ws.Range('E24').Value = '  +3.07%  '

# Row 25: LidoDAOToken
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.119'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +9.24%  '

# Row 26: Toncoin
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.336'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.18%  '

# Row 27: EthereumClassic
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.79'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.89%  '

# Row 28: Monero
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '161.79'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.34%  '

# Row 29: BitcoinCash
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '136.54'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.95%  '

# Row 30: HuobiToken
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.195'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.63%  '

# Row 31: Filecoin
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.523'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +9.74%  '

# Row 32: WrappedliquidstakedEther2.0
$ws.Range('D32').Value = '1.880.99'
$ws.Range('E32').Value = '  +1.74%  '

# Row 33: ImmutableX
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.077'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.47%  '

# Row 34: Hedera
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08567'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.32%  '

# Row 35: InternetComputer(DFINITY)
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.143'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.59%  '

# Row 36: FraxShare
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '11.51'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +10.22%  '

# Row 37: Algorand
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2732'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.29%  '

# Row 38: WEMIXTOKEN
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.931'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.04%  '

# Row 39: Aptos
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.40'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.82%  '

# Row 40: Stellar
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09122'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.72%  '

# Row 41: VeChain
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.02726'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.72%  '

# Row 42: TrustWalletToken
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.482'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.90%  '

# Row 43: TheSandbox
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7646'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.15%  '

# Row 44: Decentraland
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7152'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.61%  '

# Row 45: NEARProtocol
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.58'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.68%  '

# Row 46: EnergySwap
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.560'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.16%  '

# Row 47: PancakeSwap
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.218'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.79%  '

# Row 48: Frax
$ws.Range('E48').Value = '  +0.13%  '

# Row 49: Quant
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '140.79'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.63%  '

# Row 50: Flow
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.312'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.42%  '

# Row 51: Cronos
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07980'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.23%  '
